# Generate Report for Handback
# The c2b6a63f-fef9-4b2e-9b7c-e407c6336a56.md file has been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet status and each locale sheet's status / handback datetime.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the c2b6a63f-...md file
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# zh-cn sheet: row 3 is the c2b6a63f-...md file
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-03-18 00:38:02"

# de-de sheet: row 3 is the c2b6a63f-...md file
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-03-18 00:38:08"
